$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C holds a "Förändrad" (changed) date that was bumped by one day
# (2023-09-09 -> 2023-09-10, serial 45178 -> 45179) for every data row
# (rows 2 through 39).
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 3).Value = 45179
}
